$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 23.02.2022 21:15"

$ws.Range("D10").Value = -0.3

$ws.Range("E10").Value = 44615.87533564815
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
